$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Penk"
$ws.Range("C2").Value = "Oprm1"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.078061
$ws.Range("H2").Value = 0.234183
$ws.Range("I2").Value = 0.001800381391819829
$ws.Range("J2").Value = 0.001800381391819829
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.8377936666666667
$ws.Range("N2").Value = 2.513381
$ws.Range("O2").Value = 0.7130909380817101
$ws.Range("P2").Value = 0.7130909380817101
$ws.Range("Q2").Value = 0.06539901141366668
$ws.Range("R2").Value = 0.588591102723
$ws.Range("S2").Value = 0.001283835655597657
$ws.Range("T2").Value = 0.001283835655597657

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Penk"
$ws.Range("C3").Value = "Oprm1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.078061
$ws.Range("H3").Value = 0.234183
$ws.Range("I3").Value = 0.001800381391819829
$ws.Range("J3").Value = 0.001800381391819829
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3370826666666667
$ws.Range("N3").Value = 1.011248
$ws.Range("O3").Value = 0.2869090619182899
$ws.Range("P3").Value = 0.2869090619182899
$ws.Range("Q3").Value = 0.02631301004266667
$ws.Range("R3").Value = 0.236817090384
$ws.Range("S3").Value = 0.0005165457362221723
$ws.Range("T3").Value = 0.0005165457362221722

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Penk"
$ws.Range("C4").Value = "Oprm1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 43.033198
$ws.Range("H4").Value = 129.099594
$ws.Range("I4").Value = 0.9925080246179051
$ws.Range("J4").Value = 0.9925080246179051
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.8377936666666667
$ws.Range("N4").Value = 2.513381
$ws.Range("O4").Value = 0.7130909380817101
$ws.Range("P4").Value = 0.7130909380817101
$ws.Range("Q4").Value = 36.05294074081267
$ws.Range("R4").Value = 324.476466667314
$ws.Range("S4").Value = 0.7077484783284069
$ws.Range("T4").Value = 0.7077484783284069

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Penk"
$ws.Range("C5").Value = "Oprm1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 43.033198
$ws.Range("H5").Value = 129.099594
$ws.Range("I5").Value = 0.9925080246179051
$ws.Range("J5").Value = 0.9925080246179051
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3370826666666667
$ws.Range("N5").Value = 1.011248
$ws.Range("O5").Value = 0.2869090619182899
$ws.Range("P5").Value = 0.2869090619182899
$ws.Range("Q5").Value = 14.50574513703467
$ws.Range("R5").Value = 130.551706233312
$ws.Range("S5").Value = 0.2847595462894981
$ws.Range("T5").Value = 0.2847595462894981

$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Penk"
$ws.Range("C6").Value = "Oprm1"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.05540800000000001
$ws.Range("H6").Value = 0.166224
$ws.Range("I6").Value = 0.001277917681786719
$ws.Range("J6").Value = 0.001277917681786719
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.8377936666666667
$ws.Range("N6").Value = 2.513381
$ws.Range("O6").Value = 0.7130909380817101
$ws.Range("P6").Value = 0.7130909380817101
$ws.Range("Q6").Value = 0.04642047148266667
$ws.Range("R6").Value = 0.417784243344
$ws.Range("S6").Value = 0.000911271518496496
$ws.Range("T6").Value = 0.000911271518496496

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Penk"
$ws.Range("C7").Value = "Oprm1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.05540800000000001
$ws.Range("H7").Value = 0.166224
$ws.Range("I7").Value = 0.001277917681786719
$ws.Range("J7").Value = 0.001277917681786719
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3370826666666667
$ws.Range("N7").Value = 1.011248
$ws.Range("O7").Value = 0.2869090619182899
$ws.Range("P7").Value = 0.2869090619182899
$ws.Range("Q7").Value = 0.01867707639466667
$ws.Range("R7").Value = 0.168093687552
$ws.Range("S7").Value = 0.0003666461632902233
$ws.Range("T7").Value = 0.0003666461632902233

$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Penk"
$ws.Range("C8").Value = "Oprm1"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1913683333333333
$ws.Range("H8").Value = 0.574105
$ws.Range("I8").Value = 0.004413676308488332
$ws.Range("J8").Value = 0.004413676308488332
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.8377936666666667
$ws.Range("N8").Value = 2.513381
$ws.Range("O8").Value = 0.7130909380817101
$ws.Range("P8").Value = 0.7130909380817101
$ws.Range("Q8").Value = 0.1603271776672222
$ws.Range("R8").Value = 1.442944599005
$ws.Range("S8").Value = 0.003147352579208964
$ws.Range("T8").Value = 0.003147352579208964

$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Penk"
$ws.Range("C9").Value = "Oprm1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1913683333333333
$ws.Range("H9").Value = 0.574105
$ws.Range("I9").Value = 0.004413676308488332
$ws.Range("J9").Value = 0.004413676308488332
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.3370826666666667
$ws.Range("N9").Value = 1.011248
$ws.Range("O9").Value = 0.2869090619182899
$ws.Range("P9").Value = 0.2869090619182899
$ws.Range("Q9").Value = 0.06450694811555556
$ws.Range("R9").Value = 0.58056253304
$ws.Range("S9").Value = 0.001266323729279368
$ws.Range("T9").Value = 0.001266323729279368

